$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "U"   # Column D -> "U"
    $ws.Cells.Item($r, 6).ClearContents()  # Column F
    $ws.Cells.Item($r, 7).ClearContents()  # Column G
}

$ws.Range("D7").Select()
